$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B,C,D,E,F,I,J,K,L,M,N (column H is intentionally skipped/empty)
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# New values for rows 2..25, one row per array, matching $cols order
$data = @(
    @(1.02,1.040611444595344,1.048372287641586,0.992614727750844,1.056953125393542,1.041075353867378,1.045697114550355,1.051132316400556,0.9955398523335997,1.059689423953979,1.019074468538932),
    @(1.02,1.041853218005208,1.049355389508184,0.9936372048519299,1.058070536565577,1.04139478301384,1.046583118648308,1.051927007648588,0.9963617723202687,1.060619803631208,1.019373302413491),
    @(1.02,1.042656007461104,1.049990630505686,0.9942998659930998,1.058792967129345,1.041599446241812,1.047155210044253,1.052439721703338,0.9968940712668347,1.061220625156442,1.019566126708456),
    @(1.02,1.04299333031914,1.05025747370699,0.994578699834602,1.059096532705145,1.041685001538242,1.047395428662709,1.052654907866476,0.9971179600053012,1.061472925574931,1.019647061012254),
    @(1.02,1.04304995843282,1.050302265455106,0.994625531979634,1.059147494289601,1.041699338217479,1.047435745550194,1.05269101756632,0.9971555583673455,1.061515271257755,1.019660642687129),
    @(1.02,1.042660515449513,1.049994196909412,0.994303590798249,1.058797023952186,1.041600591341387,1.047158420989148,1.052442598439714,0.9968970624462089,1.061223997524858,1.019567208662919),
    @(1.02,1.04103125868923,1.048704716958202,0.9929600610674297,1.057330886837979,1.041183726864207,1.045996795849298,1.0514011983026,0.9958175282591056,1.060004098712096,1.019175573229331),
    @(1.02,1.038154653381421,1.046425598112071,0.9906006454969559,1.054742627298149,1.040433594004858,1.043940501983215,1.049554537158626,0.9939188001724441,1.057845252439275,1.018481292342197),
    @(1.02,1.036232942535482,1.044901463436782,0.989033133672735,1.053013829360002,1.039923002022662,1.04256323640242,1.048315556722724,0.9926553831429383,1.056399720936052,1.018015600430377),
    @(1.02,1.035399837619441,1.04424035532634,0.988355674866747,1.052264435471126,1.039699409516489,1.041965321478308,1.047777176942053,0.9921088820399291,1.055772273083646,1.017813269506942),
    @(1.02,1.035090232825162,1.043994615976423,0.9881042295826724,1.051985952833019,1.039615980398349,1.04174299414457,1.04757691245022,0.9919059725120875,1.055538980130825,1.017738011408986),
    @(1.02,1.035156651064752,1.044047335822552,0.9881581567098651,1.052045693968922,1.039633893297537,1.041790694788047,1.047619882847846,0.9919494934313052,1.055589032742109,1.017754159216645),
    @(1.02,1.035374248723011,1.044220046003161,0.9883348863814464,1.052241418559754,1.039692520933493,1.041946948640841,1.047760628867306,0.9920921077337197,1.055752993737227,1.017807050762353),
    @(1.02,1.0355082975021,1.044326435282974,0.9884438009545853,1.052361994390826,1.039728593384555,1.042043190574401,1.047847309145717,0.9921799884222134,1.055853984861977,1.017839625267878),
    @(1.02,1.036288211771243,1.044945314659187,0.9890781214508737,1.053063546845101,1.039937788314125,1.042602885200653,1.048351247172645,0.9926916645766087,1.056441330318044,1.018029014009582),
    @(1.02,1.036777163427908,1.04533321270417,0.9894763578477731,1.053503392932459,1.040068340070353,1.042953550324097,1.048666846047168,0.9930127773692701,1.056809347371073,1.018147629121247),
    @(1.02,1.03706226524155,1.045559356463389,0.9897087662937551,1.053759869427169,1.040144247277384,1.043157937975584,1.048850747004168,0.9932001317071766,1.057023858667881,1.018216749435254),
    @(1.02,1.037159461515822,1.045636446937759,0.9897880325774039,1.053847308072174,1.040170088750169,1.043227603612874,1.048913421573831,0.993264023964098,1.05709697663385,1.018240306524725),
    @(1.02,1.036724713439026,1.045291606370595,0.9894336180360677,1.053456209724805,1.040054358080933,1.042915942743674,1.048633004168692,0.9929783193494215,1.056769877826632,1.018134909661611),
    @(1.02,1.035310175875699,1.044169191988473,0.9882828385668249,1.052183785997682,1.039675266967583,1.041900942284129,1.047719190590027,0.9920501090198102,1.055704717700802,1.017791478380252),
    @(1.02,1.03441991566121,1.043462475745877,0.9875604150241495,1.051383041943929,1.039434736230507,1.041261409810462,1.047142982522483,0.9914670000341481,1.055033672408458,1.017574950987251),
    @(1.02,1.03489194446752,1.04383721576006,0.9879432794636464,1.051807600589188,1.039562453167639,1.041600567879359,1.047448599028693,0.9917760702887611,1.055389533618569,1.017689793233311),
    @(1.02,1.036748413626025,1.045310406823599,0.9894529299347244,1.053477530027275,1.04006067668607,1.042932936452,1.048648296420843,0.9929938892766441,1.056787712868481,1.018140657241454),
    @(1.02,1.038899012175266,1.047015630648893,0.9912096547607051,1.055412326783514,1.040629369745052,1.044473223794716,1.050033324592629,0.9944092447426416,1.058404469394271,1.019074468538932)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $row).Value = $rowValues[$j]
    }
}
